$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 667.8570999999999
$ws.Cells.Item(17, 10).Value = 703.8461
$ws.Cells.Item(17, 12).Value = 2111.5383
$ws.Cells.Item(17, 14).Value = -2447.5383

$ws.Cells.Item(62, 8).Value = 60353424
$ws.Cells.Item(62, 9).Value = 23819014
$ws.Cells.Item(62, 11).Value = 23819014
$ws.Cells.Item(62, 13).Value = -23818390

$ws.Cells.Item(65, 8).Value = 60353424
$ws.Cells.Item(65, 9).Value = 23819014
$ws.Cells.Item(65, 11).Value = 119095070
$ws.Cells.Item(65, 13).Value = -119091950

$ws.Cells.Item(103, 8).Value = 35716000
$ws.Cells.Item(103, 9).Value = 83333950
$ws.Cells.Item(103, 10).Value = 2530
$ws.Cells.Item(103, 11).Value = 250001850
$ws.Cells.Item(103, 12).Value = 7590
$ws.Cells.Item(103, 13).Value = -250001264
$ws.Cells.Item(103, 14).Value = -8762

$ws.Cells.Item(112, 8).Value = 564679.2
$ws.Cells.Item(112, 10).Value = 564679.2
$ws.Cells.Item(112, 12).Value = 1694037.6
$ws.Cells.Item(112, 14).Value = -1696253.6

$ws.Cells.Item(127, 8).Value = 7766.8237
$ws.Cells.Item(127, 9).Value = 17172.5
$ws.Cells.Item(127, 10).Value = 2636.4546
$ws.Cells.Item(127, 11).Value = 51517.5
$ws.Cells.Item(127, 12).Value = 7909.3638
$ws.Cells.Item(127, 13).Value = -46557.5
$ws.Cells.Item(127, 14).Value = -17829.3638

$ws.Cells.Item(135, 8).Value = 565.0270400000001
$ws.Cells.Item(135, 9).Value = 511.94116
$ws.Cells.Item(135, 10).Value = 1166.6666
$ws.Cells.Item(135, 11).Value = 4607.47044
$ws.Cells.Item(135, 12).Value = 10499.9994
$ws.Cells.Item(135, 13).Value = -2072.47044
$ws.Cells.Item(135, 14).Value = -15569.9994

$ws.Cells.Item(137, 8).Value = 9140660
$ws.Cells.Item(137, 9).Value = 778.6667
$ws.Cells.Item(137, 10).Value = 40477396
$ws.Cells.Item(137, 11).Value = 2336.0001
$ws.Cells.Item(137, 12).Value = 121432188
$ws.Cells.Item(137, 13).Value = 213.9998999999998
$ws.Cells.Item(137, 14).Value = -121437288

$ws.Cells.Item(141, 8).Value = 1445.3214
$ws.Cells.Item(141, 9).Value = 832.76
$ws.Cells.Item(141, 10).Value = 6550
$ws.Cells.Item(141, 11).Value = 2498.28
$ws.Cells.Item(141, 12).Value = 19650
$ws.Cells.Item(141, 13).Value = 2681.72
$ws.Cells.Item(141, 14).Value = -30010

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 45099532
$ws.Cells.Item(74, 9).Value = 37037720
$ws.Cells.Item(74, 11).Value = 37037720
$ws.Cells.Item(74, 13).Value = -37036846

$ws.Cells.Item(77, 8).Value = 45099532
$ws.Cells.Item(77, 9).Value = 37037720
$ws.Cells.Item(77, 11).Value = 185188600
$ws.Cells.Item(77, 13).Value = -185184232

$ws.Cells.Item(80, 8).Value = 24286
$ws.Cells.Item(80, 10).Value = 24286
$ws.Cells.Item(80, 12).Value = 24286
$ws.Cells.Item(80, 14).Value = -26282

$ws.Cells.Item(83, 8).Value = 24286
$ws.Cells.Item(83, 10).Value = 24286
$ws.Cells.Item(83, 12).Value = 72858
$ws.Cells.Item(83, 14).Value = -82842

$ws.Cells.Item(101, 8).Value = 28800
$ws.Cells.Item(101, 10).Value = 28800
$ws.Cells.Item(101, 12).Value = 28800
$ws.Cells.Item(101, 14).Value = -35290

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(63, 8).Value = 10000
$ws.Cells.Item(63, 9).Value = 10000
$ws.Cells.Item(63, 11).Value = 10000
$ws.Cells.Item(63, 13).Value = -9314

$ws.Cells.Item(66, 8).Value = 10000
$ws.Cells.Item(66, 9).Value = 10000
$ws.Cells.Item(66, 11).Value = 30000
$ws.Cells.Item(66, 13).Value = -26568

$ws.Cells.Item(134, 8).Value = 11398934
$ws.Cells.Item(134, 9).Value = 12195863
$ws.Cells.Item(134, 10).Value = 5953251
$ws.Cells.Item(134, 11).Value = 36587589
$ws.Cells.Item(134, 12).Value = 17859753
$ws.Cells.Item(134, 13).Value = -36585054
$ws.Cells.Item(134, 14).Value = -17864823

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(3, 8).Value = 660
$ws.Cells.Item(3, 9).Value = 660
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 660
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 14).Value = -547

$ws.Cells.Item(31, 8).Value = 1456611.2
$ws.Cells.Item(31, 9).Value = 1054.8438
$ws.Cells.Item(31, 10).Value = 5690957
$ws.Cells.Item(31, 11).Value = 1054.8438
$ws.Cells.Item(31, 12).Value = 5690957
$ws.Cells.Item(31, 13).Value = -759.8438000000001
$ws.Cells.Item(31, 14).Value = -5691547

$ws.Cells.Item(34, 8).Value = 1456611.2
$ws.Cells.Item(34, 9).Value = 1054.8438
$ws.Cells.Item(34, 10).Value = 5690957
$ws.Cells.Item(34, 11).Value = 1054.8438
$ws.Cells.Item(34, 12).Value = 5690957
$ws.Cells.Item(34, 13).Value = -852.8438000000001
$ws.Cells.Item(34, 14).Value = -5691361

$ws.Cells.Item(39, 8).Value = 33309.855
$ws.Cells.Item(39, 10).Value = 54792.75
$ws.Cells.Item(39, 12).Value = 54792.75
$ws.Cells.Item(39, 14).Value = -55574.75

$ws.Cells.Item(49, 8).Value = 33309.855
$ws.Cells.Item(49, 10).Value = 54792.75
$ws.Cells.Item(49, 12).Value = 54792.75
$ws.Cells.Item(49, 14).Value = -55156.75

$ws.Cells.Item(99, 8).Value = 9950
$ws.Cells.Item(99, 9).Value = 5586.8423
$ws.Cells.Item(99, 10).Value = 23766.666
$ws.Cells.Item(99, 11).Value = 5586.8423
$ws.Cells.Item(99, 12).Value = 23766.666
$ws.Cells.Item(99, 13).Value = -4088.8423
$ws.Cells.Item(99, 14).Value = -26762.666

$ws.Cells.Item(126, 8).Value = 9950
$ws.Cells.Item(126, 9).Value = 5586.8423
$ws.Cells.Item(126, 10).Value = 23766.666
$ws.Cells.Item(126, 11).Value = 16760.5269
$ws.Cells.Item(126, 12).Value = 71299.99800000001
$ws.Cells.Item(126, 13).Value = -14290.5269
$ws.Cells.Item(126, 14).Value = -76239.99800000001

$ws.Cells.Item(134, 8).Value = 678780.5600000001
$ws.Cells.Item(134, 9).Value = 817.43396
$ws.Cells.Item(134, 10).Value = 6667454.5
$ws.Cells.Item(134, 11).Value = 2452.30188
$ws.Cells.Item(134, 12).Value = 20002363.5
$ws.Cells.Item(134, 13).Value = 82.69812000000002
$ws.Cells.Item(134, 14).Value = -20007433.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 2633.2222
$ws.Cells.Item(34, 10).Value = 2924.875
$ws.Cells.Item(34, 12).Value = 8774.625
$ws.Cells.Item(34, 14).Value = -8942.625

$ws.Cells.Item(39, 8).Value = 1987.6666
$ws.Cells.Item(39, 10).Value = 2501.3333
$ws.Cells.Item(39, 12).Value = 7503.999899999999
$ws.Cells.Item(39, 14).Value = -8091.999899999999

$ws.Cells.Item(55, 8).Value = 2336.9092
$ws.Cells.Item(55, 10).Value = 2570
$ws.Cells.Item(55, 12).Value = 7710
$ws.Cells.Item(55, 14).Value = -8064

$ws.Cells.Item(112, 8).Value = 15382.125
$ws.Cells.Item(112, 9).Value = 38009
$ws.Cells.Item(112, 10).Value = 1806
$ws.Cells.Item(112, 11).Value = 114027
$ws.Cells.Item(112, 12).Value = 5418
$ws.Cells.Item(112, 13).Value = -112919
$ws.Cells.Item(112, 14).Value = -7634

$ws.Cells.Item(113, 8).Value = 1139.1323
$ws.Cells.Item(113, 9).Value = 965.1739
$ws.Cells.Item(113, 10).Value = 1228.0444
$ws.Cells.Item(113, 11).Value = 2895.5217
$ws.Cells.Item(113, 12).Value = 3684.1332
$ws.Cells.Item(113, 13).Value = -725.5217000000002
$ws.Cells.Item(113, 14).Value = -8024.1332

$ws.Cells.Item(134, 8).Value = 1828.75
$ws.Cells.Item(134, 9).Value = 2138.3333
$ws.Cells.Item(134, 10).Value = 900
$ws.Cells.Item(134, 11).Value = 6414.999899999999
$ws.Cells.Item(134, 12).Value = 2700
$ws.Cells.Item(134, 13).Value = -1344.999899999999
$ws.Cells.Item(134, 14).Value = -12840

$ws.Cells.Item(139, 8).Value = 102106
$ws.Cells.Item(139, 9).Value = 102106
$ws.Cells.Item(139, 11).Value = 306318
$ws.Cells.Item(139, 13).Value = -301178

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(41, 8).Value = 1586.8
$ws.Cells.Item(41, 9).Value = 722
$ws.Cells.Item(41, 10).Value = 2884
$ws.Cells.Item(41, 11).Value = 722
$ws.Cells.Item(41, 12).Value = 2884
$ws.Cells.Item(41, 13).Value = -367
$ws.Cells.Item(41, 14).Value = -3594

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(32, 8).Value = 3764.3333
$ws.Cells.Item(32, 9).Value = 3764.3333
$ws.Cells.Item(32, 11).Value = 3764.3333
$ws.Cells.Item(32, 13).Value = -3447.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(12, 8).Value = 36755.25
$ws.Cells.Item(12, 9).Value = 2000
$ws.Cells.Item(12, 10).Value = 48340.332
$ws.Cells.Item(12, 11).Value = 2000
$ws.Cells.Item(12, 12).Value = 48340.332
$ws.Cells.Item(12, 13).Value = -1858
$ws.Cells.Item(12, 14).Value = -48624.332

$ws.Cells.Item(14, 8).Value = 2250
$ws.Cells.Item(14, 9).Value = 1000
$ws.Cells.Item(14, 11).Value = 1000
$ws.Cells.Item(14, 13).Value = -832

$ws.Cells.Item(33, 8).Value = 7980
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 11).Value = 0
$ws.Cells.Item(33, 13).ClearContents()

$ws.Cells.Item(36, 8).Value = 7980
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 11).Value = 0
$ws.Cells.Item(36, 13).ClearContents()

$ws.Cells.Item(37, 8).Value = 8061.8
$ws.Cells.Item(37, 9).Value = 7489
$ws.Cells.Item(37, 10).Value = 8205
$ws.Cells.Item(37, 11).Value = 7489
$ws.Cells.Item(37, 12).Value = 8205
$ws.Cells.Item(37, 13).Value = -7286
$ws.Cells.Item(37, 14).Value = -8611

$ws.Cells.Item(62, 8).Value = 22740614
$ws.Cells.Item(62, 9).Value = 41686292
$ws.Cells.Item(62, 10).Value = 5799.7
$ws.Cells.Item(62, 11).Value = 41686292
$ws.Cells.Item(62, 12).Value = 5799.7
$ws.Cells.Item(62, 13).Value = -41685668
$ws.Cells.Item(62, 14).Value = -7047.7

$ws.Cells.Item(65, 8).Value = 22740614
$ws.Cells.Item(65, 9).Value = 41686292
$ws.Cells.Item(65, 10).Value = 5799.7
$ws.Cells.Item(65, 11).Value = 208431460
$ws.Cells.Item(65, 12).Value = 28998.5
$ws.Cells.Item(65, 13).Value = -208428340
$ws.Cells.Item(65, 14).Value = -35238.5

$ws.Cells.Item(81, 8).Value = 12499.619
$ws.Cells.Item(81, 9).Value = 556.5
$ws.Cells.Item(81, 10).Value = 36385.855
$ws.Cells.Item(81, 11).Value = 1113
$ws.Cells.Item(81, 12).Value = 72771.71000000001
$ws.Cells.Item(81, 13).Value = -52
$ws.Cells.Item(81, 14).Value = -74893.71000000001

$ws.Cells.Item(84, 8).Value = 12499.619
$ws.Cells.Item(84, 9).Value = 556.5
$ws.Cells.Item(84, 10).Value = 36385.855
$ws.Cells.Item(84, 11).Value = 5565
$ws.Cells.Item(84, 12).Value = 363858.55
$ws.Cells.Item(84, 13).Value = -261
$ws.Cells.Item(84, 14).Value = -374466.55
